# Scheduled market-data refresh: update Leve profit calculations (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 290
$ws.Range("I2").Value = 290
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 290
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -177
$ws.Range("N2").ClearContents()
$ws.Range("H18").Value = 2059.4375
$ws.Range("I18").Value = 925.0714
$ws.Range("K18").Value = 925.0714
$ws.Range("M18").Value = -641.0714
$ws.Range("H33").Value = 13198081
$ws.Range("I33").Value = 54095.855
$ws.Range("J33").Value = 50001240
$ws.Range("K33").Value = 54095.855
$ws.Range("L33").Value = 50001240
$ws.Range("M33").Value = -53866.855
$ws.Range("N33").Value = -50001698
$ws.Range("H55").Value = 507.65
$ws.Range("J55").Value = 574.38464
$ws.Range("L55").Value = 574.38464
$ws.Range("N55").Value = -1002.38464
$ws.Range("H62").Value = 26584.592
$ws.Range("I62").Value = 2451
$ws.Range("K62").Value = 2451
$ws.Range("M62").Value = -1827
$ws.Range("H65").Value = 26584.592
$ws.Range("I65").Value = 2451
$ws.Range("K65").Value = 12255
$ws.Range("M65").Value = -9135
$ws.Range("H80").Value = 2856189.2
$ws.Range("I80").Value = 5707941.5
$ws.Range("J80").Value = 4437
$ws.Range("K80").Value = 17123824.5
$ws.Range("L80").Value = 13311
$ws.Range("M80").Value = -17122826.5
$ws.Range("N80").Value = -15307
$ws.Range("H83").Value = 2856189.2
$ws.Range("I83").Value = 5707941.5
$ws.Range("J83").Value = 4437
$ws.Range("K83").Value = 51371473.5
$ws.Range("L83").Value = 39933
$ws.Range("M83").Value = -51366481.5
$ws.Range("N83").Value = -49917
$ws.Range("H99").Value = 360.83334
$ws.Range("I99").Value = 360.83334
$ws.Range("K99").Value = 1082.50002
$ws.Range("M99").Value = 415.4999800000001
$ws.Range("H100").Value = 6122.727
$ws.Range("I100").Value = 8192.857
$ws.Range("K100").Value = 8192.857
$ws.Range("M100").Value = -7651.857
$ws.Range("H105").Value = 26251.857
$ws.Range("J105").Value = 26251.857
$ws.Range("L105").Value = 26251.857
$ws.Range("N105").Value = -33239.857
$ws.Range("H129").Value = 1144.0869
$ws.Range("I129").Value = 938
$ws.Range("J129").Value = 2123
$ws.Range("K129").Value = 2814
$ws.Range("L129").Value = 6369
$ws.Range("M129").Value = 2186
$ws.Range("N129").Value = -16369
$ws.Range("H132").Value = 5052.7026
$ws.Range("J132").Value = 5759.125
$ws.Range("L132").Value = 17277.375
$ws.Range("N132").Value = -22337.375
$ws.Range("H138").Value = 2294.6562
$ws.Range("J138").Value = 2824.842
$ws.Range("L138").Value = 8474.526
$ws.Range("N138").Value = -18754.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1674.4286
$ws.Range("I2").Value = 1419.6364
$ws.Range("J2").Value = 2608.6667
$ws.Range("K2").Value = 1419.6364
$ws.Range("L2").Value = 2608.6667
$ws.Range("M2").Value = -1306.6364
$ws.Range("N2").Value = -2834.6667
$ws.Range("H32").Value = 16509.217
$ws.Range("I32").Value = 9263.134
$ws.Range("K32").Value = 9263.134
$ws.Range("M32").Value = -8976.134
$ws.Range("H61").Value = 4130.923
$ws.Range("I61").Value = 2712.5715
$ws.Range("J61").Value = 5785.6665
$ws.Range("K61").Value = 2712.5715
$ws.Range("L61").Value = 5785.6665
$ws.Range("M61").Value = -2500.5715
$ws.Range("N61").Value = -6209.6665
$ws.Range("H116").Value = 1674.4286
$ws.Range("I116").Value = 1419.6364
$ws.Range("J116").Value = 2608.6667
$ws.Range("K116").Value = 1419.6364
$ws.Range("L116").Value = 2608.6667
$ws.Range("M116").Value = 874.3635999999999
$ws.Range("N116").Value = -7196.6667
$ws.Range("H123").Value = 48301
$ws.Range("I123").Value = 43211
$ws.Range("J123").Value = 49997.668
$ws.Range("K123").Value = 43211
$ws.Range("L123").Value = 49997.668
$ws.Range("M123").Value = -38311
$ws.Range("N123").Value = -59797.668
$ws.Range("H135").Value = 103841.6
$ws.Range("I135").Value = 49211
$ws.Range("J135").Value = 117499.25
$ws.Range("K135").Value = 49211
$ws.Range("L135").Value = 117499.25
$ws.Range("M135").Value = -44141
$ws.Range("N135").Value = -127639.25
$ws.Range("H136").Value = 4130.923
$ws.Range("I136").Value = 2712.5715
$ws.Range("J136").Value = 5785.6665
$ws.Range("K136").Value = 8137.7145
$ws.Range("L136").Value = 17356.9995
$ws.Range("M136").Value = -5587.7145
$ws.Range("N136").Value = -22456.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1674.4286
$ws.Range("I3").Value = 1419.6364
$ws.Range("J3").Value = 2608.6667
$ws.Range("K3").Value = 1419.6364
$ws.Range("L3").Value = 2608.6667
$ws.Range("M3").Value = -1305.6364
$ws.Range("N3").Value = -2836.6667
$ws.Range("H94").Value = 6581084
$ws.Range("I94").Value = 8622460
$ws.Range("J94").Value = 3317.111
$ws.Range("K94").Value = 8622460
$ws.Range("L94").Value = 3317.111
$ws.Range("M94").Value = -8622009
$ws.Range("N94").Value = -4219.111
$ws.Range("H100").Value = 26826.428
$ws.Range("J100").Value = 26826.428
$ws.Range("L100").Value = 26826.428
$ws.Range("N100").Value = -28990.428
$ws.Range("H103").Value = 16261.5
$ws.Range("J103").Value = 16261.5
$ws.Range("L103").Value = 16261.5
$ws.Range("N103").Value = -18605.5
$ws.Range("H105").Value = 2233.56
$ws.Range("I105").Value = 1998.7333
$ws.Range("K105").Value = 1998.7333
$ws.Range("M105").Value = -251.7333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 20139.166
$ws.Range("J28").Value = 20139.166
$ws.Range("L28").Value = 20139.166
$ws.Range("N28").Value = -20629.166
$ws.Range("H43").Value = 14404.692
$ws.Range("J43").Value = 14404.692
$ws.Range("L43").Value = 14404.692
$ws.Range("N43").Value = -14772.692
$ws.Range("H58").Value = 419202.66
$ws.Range("I58").Value = 2159.5
$ws.Range("J58").Value = 1253289
$ws.Range("K58").Value = 2159.5
$ws.Range("L58").Value = 1253289
$ws.Range("M58").Value = -1956.5
$ws.Range("N58").Value = -1253695
$ws.Range("H62").Value = 104380
$ws.Range("I62").Value = 3450
$ws.Range("J62").Value = 171666.67
$ws.Range("K62").Value = 3450
$ws.Range("L62").Value = 171666.67
$ws.Range("M62").Value = -2826
$ws.Range("N62").Value = -172914.67
$ws.Range("H65").Value = 104380
$ws.Range("I65").Value = 3450
$ws.Range("J65").Value = 171666.67
$ws.Range("K65").Value = 17250
$ws.Range("L65").Value = 858333.3500000001
$ws.Range("M65").Value = -14130
$ws.Range("N65").Value = -864573.3500000001
$ws.Range("H94").Value = 1702.5834
$ws.Range("J94").Value = 1792.4445
$ws.Range("L94").Value = 1792.4445
$ws.Range("N94").Value = -2694.4445
$ws.Range("H101").Value = 14404.692
$ws.Range("J101").Value = 14404.692
$ws.Range("L101").Value = 14404.692
$ws.Range("N101").Value = -20894.692
$ws.Range("H104").Value = 57927
$ws.Range("J104").Value = 60285
$ws.Range("L104").Value = 60285
$ws.Range("N104").Value = -65527
$ws.Range("H105").Value = 808.2857
$ws.Range("I105").Value = 831.25
$ws.Range("K105").Value = 831.25
$ws.Range("M105").Value = 915.75
$ws.Range("H107").Value = 643.25
$ws.Range("I107").Value = 639.6667
$ws.Range("K107").Value = 639.6667
$ws.Range("M107").Value = 1280.3333
$ws.Range("H132").Value = 272167.75
$ws.Range("I132").Value = 1950.1945
$ws.Range("J132").Value = 10000000
$ws.Range("K132").Value = 5850.583500000001
$ws.Range("L132").Value = 30000000
$ws.Range("M132").Value = -3320.583500000001
$ws.Range("N132").Value = -30005060
$ws.Range("H136").Value = 419202.66
$ws.Range("I136").Value = 2159.5
$ws.Range("J136").Value = 1253289
$ws.Range("K136").Value = 6478.5
$ws.Range("L136").Value = 3759867
$ws.Range("M136").Value = -3928.5
$ws.Range("N136").Value = -3764967

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 908.28
$ws.Range("I2").Value = 875
$ws.Range("K2").Value = 5250
$ws.Range("M2").Value = -5137
$ws.Range("H4").Value = 45582704
$ws.Range("I4").Value = 48481470
$ws.Range("K4").Value = 145444410
$ws.Range("M4").Value = -145444298
$ws.Range("H33").Value = 2292.3809
$ws.Range("I33").Value = 1861.5834
$ws.Range("K33").Value = 11169.5004
$ws.Range("M33").Value = -10886.5004
$ws.Range("H60").Value = 40
$ws.Range("I60").Value = 40
$ws.Range("K60").Value = 120
$ws.Range("M60").Value = 131
$ws.Range("H68").Value = 1859.2
$ws.Range("J68").Value = 2364.6667
$ws.Range("L68").Value = 7094.000100000001
$ws.Range("N68").Value = -8716.000100000001
$ws.Range("H71").Value = 1859.2
$ws.Range("J71").Value = 2364.6667
$ws.Range("L71").Value = 21282.0003
$ws.Range("N71").Value = -29394.0003
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H98").Value = 1224.6666
$ws.Range("I98").Value = 1774.5
$ws.Range("K98").Value = 5323.5
$ws.Range("M98").Value = -3825.5
$ws.Range("H113").Value = 676.375
$ws.Range("I113").Value = 487.25
$ws.Range("J113").Value = 865.5
$ws.Range("K113").Value = 1461.75
$ws.Range("L113").Value = 2596.5
$ws.Range("M113").Value = 708.25
$ws.Range("N113").Value = -6936.5
$ws.Range("H120").Value = 15000
$ws.Range("J120").Value = 16666.666
$ws.Range("L120").Value = 49999.99800000001
$ws.Range("N120").Value = -59675.99800000001
$ws.Range("H137").Value = 1985.1111
$ws.Range("I137").Value = 914.9
$ws.Range("J137").Value = 3322.875
$ws.Range("K137").Value = 2744.7
$ws.Range("L137").Value = 9968.625
$ws.Range("M137").Value = 2355.3
$ws.Range("N137").Value = -20168.625
$ws.Range("H140").Value = 1947.3334
$ws.Range("I140").Value = 619.8
$ws.Range("K140").Value = 1859.4
$ws.Range("M140").Value = 3320.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 286.72726
$ws.Range("I2").Value = 66
$ws.Range("K2").Value = 66
$ws.Range("M2").Value = 47
$ws.Range("H125").Value = 68441.39999999999
$ws.Range("I125").Value = 48211
$ws.Range("J125").Value = 73499
$ws.Range("K125").Value = 48211
$ws.Range("L125").Value = 73499
$ws.Range("M125").Value = -45751
$ws.Range("N125").Value = -78419
$ws.Range("H132").Value = 6767.1055
$ws.Range("I132").Value = 5998.8184
$ws.Range("K132").Value = 17996.4552
$ws.Range("M132").Value = -15466.4552
$ws.Range("H134").Value = 71482.836
$ws.Range("J134").Value = 71482.836
$ws.Range("L134").Value = 214448.508
$ws.Range("N134").Value = -219518.508
$ws.Range("H136").Value = 38436.668
$ws.Range("J136").Value = 38436.668
$ws.Range("L136").Value = 115310.004
$ws.Range("N136").Value = -120410.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1466.5
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 1599.75
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 1599.75
$ws.Range("M22").Value = -905
$ws.Range("N22").Value = -2189.75
$ws.Range("H27").Value = 1466.5
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 1599.75
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 1599.75
$ws.Range("M27").Value = -1093
$ws.Range("N27").Value = -1813.75
$ws.Range("H40").Value = 3033.2
$ws.Range("I40").Value = 2999.889
$ws.Range("J40").Value = 3083.1667
$ws.Range("K40").Value = 2999.889
$ws.Range("L40").Value = 3083.1667
$ws.Range("M40").Value = -2863.889
$ws.Range("N40").Value = -3355.1667
$ws.Range("H55").Value = 507.4
$ws.Range("I55").Value = 645.3333
$ws.Range("J55").Value = 300.5
$ws.Range("K55").Value = 645.3333
$ws.Range("L55").Value = 300.5
$ws.Range("M55").Value = -472.3333
$ws.Range("N55").Value = -646.5
$ws.Range("H93").Value = 17546730
$ws.Range("I93").Value = 23811978
$ws.Range("K93").Value = 23811978
$ws.Range("M93").Value = -23810730
$ws.Range("H100").Value = 2956.7144
$ws.Range("H122").Value = 5459.7
$ws.Range("I122").Value = 3245
$ws.Range("J122").Value = 8166.5557
$ws.Range("K122").Value = 9735
$ws.Range("L122").Value = 24499.6671
$ws.Range("M122").Value = -7285
$ws.Range("N122").Value = -29399.6671
$ws.Range("H135").Value = 59999
$ws.Range("J135").Value = 59999
$ws.Range("L135").Value = 59999
$ws.Range("N135").Value = -70139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1666.6666
$ws.Range("I7").Value = 1666.6666
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1666.6666
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1553.6666
$ws.Range("N7").ClearContents()
$ws.Range("H96").Value = 2959.4
$ws.Range("J96").Value = 2949.1667
$ws.Range("L96").Value = 2949.1667
$ws.Range("N96").Value = -5695.1667
$ws.Range("H101").Value = 23683.334
$ws.Range("J101").Value = 23683.334
$ws.Range("L101").Value = 23683.334
$ws.Range("N101").Value = -30173.334
$ws.Range("H104").Value = 10102.5
$ws.Range("J104").Value = 10102.5
$ws.Range("L104").Value = 10102.5
$ws.Range("N104").Value = -17090.5
